# Update TPM-derived ligand/receptor expression & edge-weight figures on
# Sheet1 (Psen1-Ncstn LR-pair table) to reflect the new TPM inputs.
# Columns G,H (ligand avg/total expr) and I,J (ligand specificity) vary by
# "Sending cluster" (col A); columns M,N (receptor avg/total expr) and O,P
# (receptor specificity) vary by "Target cluster" (col D); the edge columns
# Q=G*M, R=H*N, S=I*O, T=J*P are recomputed accordingly for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.409481333333334
$ws.Range("H2").Value = 28.228444
$ws.Range("I2").Value = 0.2433300530093958
$ws.Range("J2").Value = 0.2433300530093958
$ws.Range("M2").Value = 13.11484166666667
$ws.Range("N2").Value = 39.344525
$ws.Range("O2").Value = 0.2389043281085165
$ws.Range("P2").Value = 0.2389043281085165
$ws.Range("Q2").Value = 123.4038578521222
$ws.Range("R2").Value = 1110.6347206691
$ws.Range("S2").Value = 0.05813260282281942
$ws.Range("T2").Value = 0.05813260282281942

$ws.Range("G3").Value = 9.409481333333334
$ws.Range("H3").Value = 28.228444
$ws.Range("I3").Value = 0.2433300530093958
$ws.Range("J3").Value = 0.2433300530093958
$ws.Range("O3").Value = 0.5367702700792449
$ws.Range("P3").Value = 0.5367702700792449
$ws.Range("Q3").Value = 277.2638010895169
$ws.Range("R3").Value = 2495.374209805652
$ws.Range("S3").Value = 0.1306123382722504
$ws.Range("T3").Value = 0.1306123382722504

$ws.Range("G4").Value = 9.409481333333334
$ws.Range("H4").Value = 28.228444
$ws.Range("I4").Value = 0.2433300530093958
$ws.Range("J4").Value = 0.2433300530093958
$ws.Range("O4").Value = 0.2243254018122386
$ws.Range("P4").Value = 0.2243254018122386
$ws.Range("Q4").Value = 115.87324606896
$ws.Range("R4").Value = 1042.85921462064
$ws.Range("S4").Value = 0.05458511191432603
$ws.Range("T4").Value = 0.05458511191432604

$ws.Range("I5").Value = 0.5069354697952918
$ws.Range("J5").Value = 0.5069354697952919
$ws.Range("M5").Value = 13.11484166666667
$ws.Range("N5").Value = 39.344525
$ws.Range("O5").Value = 0.2389043281085165
$ws.Range("P5").Value = 0.2389043281085165
$ws.Range("Q5").Value = 257.0902849078056
$ws.Range("R5").Value = 2313.81256417025
$ws.Range("S5").Value = 0.1211090778058194
$ws.Range("T5").Value = 0.1211090778058194

$ws.Range("I6").Value = 0.5069354697952918
$ws.Range("J6").Value = 0.5069354697952919
$ws.Range("O6").Value = 0.5367702700792449
$ws.Range("P6").Value = 0.5367702700792449
$ws.Range("S6").Value = 0.2721078890347677
$ws.Range("T6").Value = 0.2721078890347677

$ws.Range("I7").Value = 0.5069354697952918
$ws.Range("J7").Value = 0.5069354697952919
$ws.Range("O7").Value = 0.2243254018122386
$ws.Range("P7").Value = 0.2243254018122386
$ws.Range("S7").Value = 0.1137185029547048
$ws.Range("T7").Value = 0.1137185029547048

$ws.Range("G8").Value = 9.657138
$ws.Range("I8").Value = 0.2497344771953123
$ws.Range("J8").Value = 0.2497344771953124
$ws.Range("M8").Value = 13.11484166666667
$ws.Range("N8").Value = 39.344525
$ws.Range("O8").Value = 0.2389043281085165
$ws.Range("P8").Value = 0.2389043281085165
$ws.Range("Q8").Value = 126.65183582315
$ws.Range("R8").Value = 1139.86652240835
$ws.Range("S8").Value = 0.05966264747987773
$ws.Range("T8").Value = 0.05966264747987774

$ws.Range("G9").Value = 9.657138
$ws.Range("I9").Value = 0.2497344771953123
$ws.Range("J9").Value = 0.2497344771953124
$ws.Range("O9").Value = 0.5367702700792449
$ws.Range("P9").Value = 0.5367702700792449
$ws.Range("Q9").Value = 284.5613583440179
$ws.Range("S9").Value = 0.1340500427722268
$ws.Range("T9").Value = 0.1340500427722268

$ws.Range("G10").Value = 9.657138
$ws.Range("I10").Value = 0.2497344771953123
$ws.Range("J10").Value = 0.2497344771953124
$ws.Range("O10").Value = 0.2243254018122386
$ws.Range("P10").Value = 0.2243254018122386
$ws.Range("S10").Value = 0.05602178694320777
$ws.Range("T10").Value = 0.05602178694320779
